$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2365.875
$ws.Range("I19").Value = 2116.6667
$ws.Range("J19").Value = 2515.4
$ws.Range("K19").Value = 2116.6667
$ws.Range("L19").Value = 2515.4
$ws.Range("M19").Value = -1941.6667
$ws.Range("N19").Value = -2865.4
$ws.Range("H62").Value = 9808941
$ws.Range("I62").Value = 10421744
$ws.Range("K62").Value = 10421744
$ws.Range("M62").Value = -10421120
$ws.Range("H64").Value = 166671620
$ws.Range("I64").Value = 333333340
$ws.Range("J64").Value = 9900
$ws.Range("K64").Value = 333333340
$ws.Range("L64").Value = 9900
$ws.Range("M64").Value = -333333092
$ws.Range("N64").Value = -10396
$ws.Range("H65").Value = 9808941
$ws.Range("I65").Value = 10421744
$ws.Range("K65").Value = 52108720
$ws.Range("M65").Value = -52105600
$ws.Range("H67").Value = 166671620
$ws.Range("I67").Value = 333333340
$ws.Range("J67").Value = 9900
$ws.Range("K67").Value = 333333340
$ws.Range("L67").Value = 9900
$ws.Range("M67").Value = -333332482
$ws.Range("N67").Value = -11616
$ws.Range("H70").Value = 14500.5
$ws.Range("J70").Value = 14500.5
$ws.Range("L70").Value = 43501.5
$ws.Range("N70").Value = -44041.5
$ws.Range("H73").Value = 14500.5
$ws.Range("J73").Value = 14500.5
$ws.Range("L73").Value = 43501.5
$ws.Range("N73").Value = -45373.5
$ws.Range("H98").Value = 3303.4
$ws.Range("I98").Value = 3003.6667
$ws.Range("J98").Value = 3753
$ws.Range("K98").Value = 3003.6667
$ws.Range("L98").Value = 3753
$ws.Range("M98").Value = -1505.6667
$ws.Range("N98").Value = -6749
$ws.Range("H106").Value = 4749.6665
$ws.Range("I106").Value = 3249.7144
$ws.Range("K106").Value = 3249.7144
$ws.Range("M106").Value = -2618.7144
$ws.Range("H122").Value = 3303.4
$ws.Range("I122").Value = 3003.6667
$ws.Range("J122").Value = 3753
$ws.Range("K122").Value = 9011.000100000001
$ws.Range("L122").Value = 11259
$ws.Range("M122").Value = -6561.000100000001
$ws.Range("N122").Value = -16159
$ws.Range("H132").Value = 294537.75
$ws.Range("I132").Value = 336030.94
$ws.Range("K132").Value = 1008092.82
$ws.Range("M132").Value = -1005562.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3530.0334
$ws.Range("I2").Value = 1473.88
$ws.Range("K2").Value = 1473.88
$ws.Range("M2").Value = -1360.88
$ws.Range("H45").Value = 1038.8
$ws.Range("I45").Value = 1066.6666
$ws.Range("K45").Value = 1066.6666
$ws.Range("M45").Value = -689.6666
$ws.Range("H97").Value = 731.52
$ws.Range("I97").Value = 387.30768
$ws.Range("J97").Value = 1104.4166
$ws.Range("K97").Value = 387.30768
$ws.Range("L97").Value = 1104.4166
$ws.Range("M97").Value = 108.69232
$ws.Range("N97").Value = -2096.4166
$ws.Range("H116").Value = 3530.0334
$ws.Range("I116").Value = 1473.88
$ws.Range("K116").Value = 1473.88
$ws.Range("M116").Value = 820.1199999999999
$ws.Range("H132").Value = 1017790.1
$ws.Range("I132").Value = 1280574.1
$ws.Range("K132").Value = 3841722.3
$ws.Range("M132").Value = -3839192.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3530.0334
$ws.Range("I3").Value = 1473.88
$ws.Range("K3").Value = 1473.88
$ws.Range("M3").Value = -1359.88
$ws.Range("H86").Value = 2733.4285
$ws.Range("I86").Value = 1291.6666
$ws.Range("K86").Value = 1291.6666
$ws.Range("M86").Value = -168.6666
$ws.Range("H89").Value = 2733.4285
$ws.Range("I89").Value = 1291.6666
$ws.Range("K89").Value = 6458.333000000001
$ws.Range("M89").Value = -842.3330000000005
$ws.Range("H94").Value = 34913.72
$ws.Range("I94").Value = 1365.5
$ws.Range("K94").Value = 1365.5
$ws.Range("M94").Value = -914.5
$ws.Range("H99").Value = 8285.236000000001
$ws.Range("I99").Value = 7876.5
$ws.Range("K99").Value = 7876.5
$ws.Range("M99").Value = -6378.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 55577636
$ws.Range("I58").Value = 166683330
$ws.Range("K58").Value = 166683330
$ws.Range("M58").Value = -166683127
$ws.Range("H62").Value = 9644.777
$ws.Range("J62").Value = 4200
$ws.Range("L62").Value = 4200
$ws.Range("N62").Value = -5448
$ws.Range("H65").Value = 9644.777
$ws.Range("J65").Value = 4200
$ws.Range("L65").Value = 21000
$ws.Range("N65").Value = -27240
$ws.Range("H99").Value = 7939487.5
$ws.Range("I99").Value = 13891802
$ws.Range("J99").Value = 3068.8333
$ws.Range("K99").Value = 13891802
$ws.Range("L99").Value = 3068.8333
$ws.Range("M99").Value = -13890304
$ws.Range("N99").Value = -6064.8333
$ws.Range("H126").Value = 7939487.5
$ws.Range("I126").Value = 13891802
$ws.Range("J126").Value = 3068.8333
$ws.Range("K126").Value = 41675406
$ws.Range("L126").Value = 9206.499899999999
$ws.Range("M126").Value = -41672936
$ws.Range("N126").Value = -14146.4999
$ws.Range("H136").Value = 55577636
$ws.Range("I136").Value = 166683330
$ws.Range("K136").Value = 500049990
$ws.Range("M136").Value = -500047440
$ws.Range("H137").Value = 113995
$ws.Range("J137").Value = 113995
$ws.Range("L137").Value = 113995
$ws.Range("N137").Value = -124195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 742.58826
$ws.Range("I113").Value = 409
$ws.Range("J113").Value = 924.5454999999999
$ws.Range("K113").Value = 1227
$ws.Range("L113").Value = 2773.6365
$ws.Range("M113").Value = 943
$ws.Range("N113").Value = -7113.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4527.8125
$ws.Range("I80").Value = 2745.8462
$ws.Range("J80").Value = 12249.667
$ws.Range("K80").Value = 2745.8462
$ws.Range("L80").Value = 12249.667
$ws.Range("M80").Value = -1747.8462
$ws.Range("N80").Value = -14245.667
$ws.Range("H83").Value = 4527.8125
$ws.Range("I83").Value = 2745.8462
$ws.Range("J83").Value = 12249.667
$ws.Range("K83").Value = 13729.231
$ws.Range("L83").Value = 61248.335
$ws.Range("M83").Value = -8737.231
$ws.Range("N83").Value = -71232.33499999999
$ws.Range("H102").Value = 1005377.25
$ws.Range("I102").Value = 1670690.8
$ws.Range("K102").Value = 1670690.8
$ws.Range("M102").Value = -1669068.8
$ws.Range("H113").Value = 8591.895
$ws.Range("I113").Value = 5660.7144
$ws.Range("J113").Value = 10301.75
$ws.Range("K113").Value = 5660.7144
$ws.Range("L113").Value = 10301.75
$ws.Range("M113").Value = -3490.7144
$ws.Range("N113").Value = -14641.75
$ws.Range("H122").Value = 5289.0557
$ws.Range("J122").Value = 6241.3
$ws.Range("L122").Value = 18723.9
$ws.Range("N122").Value = -23623.9
$ws.Range("H132").Value = 62506816
$ws.Range("I132").Value = 100007210
$ws.Range("J132").Value = 6166.5
$ws.Range("K132").Value = 300021630
$ws.Range("L132").Value = 18499.5
$ws.Range("M132").Value = -300019100
$ws.Range("N132").Value = -23559.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12389
$ws.Range("I7").Value = 7917.5
$ws.Range("K7").Value = 7917.5
$ws.Range("M7").Value = -7805.5
$ws.Range("H126").Value = 12389
$ws.Range("I126").Value = 7917.5
$ws.Range("K126").Value = 23752.5
$ws.Range("M126").Value = -21282.5
$ws.Range("H136").Value = 47627790
$ws.Range("I136").Value = 142868880
$ws.Range("J136").Value = 7242.9287
$ws.Range("K136").Value = 428606640
$ws.Range("L136").Value = 21728.7861
$ws.Range("M136").Value = -428604090
$ws.Range("N136").Value = -26828.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12078.556
$ws.Range("I122").Value = 3978.4119
$ws.Range("J122").Value = 25848.8
$ws.Range("K122").Value = 11935.2357
$ws.Range("L122").Value = 77546.39999999999
$ws.Range("M122").Value = -9485.235700000001
$ws.Range("N122").Value = -82446.39999999999
$ws.Range("H136").Value = 71475420
$ws.Range("I136").Value = 166737230
$ws.Range("J136").Value = 29062.5
$ws.Range("K136").Value = 500211690
$ws.Range("L136").Value = 87187.5
$ws.Range("M136").Value = -500209140
$ws.Range("N136").Value = -92287.5
$ws.Range("H141").Value = 84166.664
$ws.Range("J141").Value = 84166.664
$ws.Range("L141").Value = 84166.664
$ws.Range("N141").Value = -94526.664
